# Rename the (only) sheet from "Sheet" to "Seed_42"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Seed_42"

# Header row (row 1) - 18 metric names across columns A..R
$headers = @("info/learner/policy4/learner_stats/total_loss","info/learner/policy2/learner_stats/total_loss","info/learner/policy1/learner_stats/total_loss","info/learner/policy3/learner_stats/total_loss","info/learner/policy0/learner_stats/total_loss","info/learner/policy4/learner_stats/entropy","info/learner/policy4/learner_stats/entropy_coeff","info/learner/policy2/learner_stats/entropy","info/learner/policy2/learner_stats/entropy_coeff","info/learner/policy1/learner_stats/entropy","info/learner/policy1/learner_stats/entropy_coeff","info/learner/policy3/learner_stats/entropy","info/learner/policy3/learner_stats/entropy_coeff","info/learner/policy0/learner_stats/entropy","info/learner/policy0/learner_stats/entropy_coeff","training_iteration","episode_len_mean","episode_reward_mean")

# Data rows 2..4 (one entry per column A..R); $null means "leave the cell empty"
$row2 = @(9.960103416442871,9.955968570709228,9.974653434753415,9.961173439025879,9.967659664154052,0.6919238865375519,0.03,0.6880520880222321,0.03,0.6918735921382904,0.03,0.6914673626422883,0.03,0.6917753994464875,0.03,1,$null,$null)
$row3 = @(9.959320354461671,9.9526198387146,9.750195407867432,9.939899158477782,9.860901165008546,0.6796949326992034,0.03,0.6419391572475434,0.03,0.6885584831237793,0.03,0.6640703916549683,0.03,0.669375067949295,0.03,2,80,-15758.57142857143)
$row4 = @(9.821837902069092,9.939698123931883,9.779951477050782,9.944232749938966,9.973601245880127,0.6467763006687164,0.03,0.6223031640052795,0.03,0.6555776000022888,0.03,0.6057752430438995,0.03,0.6262050271034241,0.03,3,80,-16992.14285714286)

$dataRows = @($row2, $row3, $row4)

# Write header values
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Apply the bold / thin-border / centered-top style to A1, then propagate the
# exact same style to the rest of the header row via copy/paste-special so we
# don't create a pile of transient intermediate cell styles (one per partial
# formatting step x per cell) in styles.xml.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160

$a1.Copy()
$ws.Range("B1:R1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the data rows (rows 2..4), skipping cells whose value is $null
for ($r = 0; $r -lt $dataRows.Length; $r++) {
    $rowValues = $dataRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $val = $rowValues[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $val
        }
    }
}

Write-Host "done"
